$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2,1).Value = 'Última actualización: 18:12:36'
$ws.Cells.Item(3,1).Value = 'Total filas: 438'
$ws.Cells.Item(45,1).Value = '05:20:00'
$ws.Cells.Item(45,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(45,4).Value = 116
$ws.Cells.Item(46,1).Value = '06:52:23'
$ws.Cells.Item(46,3).Value = '16_SANTA ANA'
$ws.Cells.Item(46,4).Value = 24
$ws.Cells.Item(73,1).Value = '07:46:15'
$ws.Cells.Item(73,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(73,4).Value = 37
$ws.Cells.Item(74,1).Value = '06:52:23'
$ws.Cells.Item(74,3).Value = '215B_EL PATO'
$ws.Cells.Item(74,4).Value = 91
$ws.Cells.Item(142,1).Value = '08:50:00'
$ws.Cells.Item(142,3).Value = '215A_EL PATO'
$ws.Cells.Item(142,4).Value = 97
$ws.Cells.Item(143,1).Value = '09:38:04'
$ws.Cells.Item(143,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(143,4).Value = 49
$ws.Cells.Item(224,3).Value = '27_EL RETIRO'
$ws.Cells.Item(225,1).Value = '11:51:05'
$ws.Cells.Item(225,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(225,4).Value = 46
$ws.Cells.Item(226,1).Value = '10:57:58'
$ws.Cells.Item(226,3).Value = '17_179 Y 38'
$ws.Cells.Item(226,4).Value = 100
$ws.Cells.Item(251,1).Value = '12:44:21'
$ws.Cells.Item(251,3).Value = '10_OLMOS'
$ws.Cells.Item(251,4).Value = 37
$ws.Cells.Item(252,1).Value = '11:51:05'
$ws.Cells.Item(252,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(252,4).Value = 90
$ws.Cells.Item(264,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(265,3).Value = '215A_EL PATO'
$ws.Cells.Item(346,1).Value = '16:14:52'
$ws.Cells.Item(346,3).Value = '16_SANTA ANA'
$ws.Cells.Item(346,4).Value = 20
$ws.Cells.Item(347,1).Value = '15:51:40'
$ws.Cells.Item(347,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(347,4).Value = 43
$ws.Cells.Item(360,1).Value = '15:19:52'
$ws.Cells.Item(360,3).Value = '17_179 Y 38'
$ws.Cells.Item(360,4).Value = 97
$ws.Cells.Item(361,1).Value = '16:14:52'
$ws.Cells.Item(361,3).Value = '10_OLMOS'
$ws.Cells.Item(361,4).Value = 42
$ws.Cells.Item(386,3).Value = '16_SANTA ANA'
$ws.Cells.Item(387,3).Value = '17_ROMERO'
$ws.Cells.Item(388,3).Value = '215B_EL PATO'
$ws.Cells.Item(396,1).Value = '16:14:52'
$ws.Cells.Item(396,3).Value = '81_EL PELIGRO'
$ws.Cells.Item(396,4).Value = 98
$ws.Cells.Item(397,1).Value = '17:39:57'
$ws.Cells.Item(397,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(397,4).Value = 13
$ws.Cells.Item(408,1).Value = '18:12:36'
$ws.Cells.Item(408,2).Value = '18:12'
$ws.Cells.Item(408,3).Value = '17_ROMERO'
$ws.Cells.Item(408,4).Value = 0
$ws.Cells.Item(410,1).Value = '17:39:57'
$ws.Cells.Item(410,2).Value = '18:15'
$ws.Cells.Item(410,4).Value = 36
$ws.Cells.Item(411,3).Value = '15_ABASTO'
$ws.Cells.Item(412,1).Value = '17:53:46'
$ws.Cells.Item(412,2).Value = '18:16'
$ws.Cells.Item(412,3).Value = '10_OLMOS'
$ws.Cells.Item(412,4).Value = 23
$ws.Cells.Item(413,1).Value = '17:14:54'
$ws.Cells.Item(413,2).Value = '18:20'
$ws.Cells.Item(413,4).Value = 66
$ws.Cells.Item(414,1).Value = '18:12:36'
$ws.Cells.Item(414,2).Value = '18:20'
$ws.Cells.Item(414,3).Value = '16_SANTA ANA'
$ws.Cells.Item(414,4).Value = 8
$ws.Cells.Item(415,2).Value = '18:21'
$ws.Cells.Item(415,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(415,4).Value = 109
$ws.Cells.Item(416,1).Value = '17:39:57'
$ws.Cells.Item(416,2).Value = '18:24'
$ws.Cells.Item(416,3).Value = '14_ABASTO'
$ws.Cells.Item(416,4).Value = 45
$ws.Cells.Item(417,1).Value = '16:32:38'
$ws.Cells.Item(417,2).Value = '18:27'
$ws.Cells.Item(417,3).Value = '215C_EL PATO'
$ws.Cells.Item(417,4).Value = 115
$ws.Cells.Item(418,2).Value = '18:28'
$ws.Cells.Item(418,3).Value = '215C_EL PATO'
$ws.Cells.Item(418,4).Value = 103
$ws.Cells.Item(419,1).Value = '18:12:36'
$ws.Cells.Item(419,2).Value = '18:30'
$ws.Cells.Item(419,3).Value = '16_SANTA ANA'
$ws.Cells.Item(419,4).Value = 18
$ws.Cells.Item(420,1).Value = '17:14:54'
$ws.Cells.Item(420,2).Value = '18:31'
$ws.Cells.Item(420,3).Value = '11X44_ETCHEVERRY'
$ws.Cells.Item(420,4).Value = 77
$ws.Cells.Item(421,1).Value = '16:45:22'
$ws.Cells.Item(421,2).Value = '18:32'
$ws.Cells.Item(421,3).Value = '11X44_ETCHEVERRY'
$ws.Cells.Item(421,4).Value = 107
$ws.Cells.Item(422,1).Value = '17:39:57'
$ws.Cells.Item(422,2).Value = '18:36'
$ws.Cells.Item(422,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(422,4).Value = 57
$ws.Cells.Item(423,1).Value = '17:53:46'
$ws.Cells.Item(423,2).Value = '18:40'
$ws.Cells.Item(423,3).Value = '15_ABASTO'
$ws.Cells.Item(423,4).Value = 47
$ws.Cells.Item(424,1).Value = '18:12:36'
$ws.Cells.Item(424,2).Value = '18:40'
$ws.Cells.Item(424,3).Value = '14_ABASTO'
$ws.Cells.Item(424,4).Value = 28
$ws.Cells.Item(425,2).Value = '18:47'
$ws.Cells.Item(425,3).Value = '14X44_ABASTO'
$ws.Cells.Item(425,4).Value = 93
$ws.Cells.Item(426,1).Value = '16:52:27'
$ws.Cells.Item(426,2).Value = '18:48'
$ws.Cells.Item(426,3).Value = '14X44_ABASTO'
$ws.Cells.Item(426,4).Value = 116
$ws.Cells.Item(427,1).Value = '18:12:36'
$ws.Cells.Item(427,2).Value = '18:52'
$ws.Cells.Item(427,3).Value = '15_ABASTO'
$ws.Cells.Item(427,4).Value = 40
$ws.Cells.Item(428,1).Value = '18:12:36'
$ws.Cells.Item(428,2).Value = '18:56'
$ws.Cells.Item(428,3).Value = '10_OLMOS'
$ws.Cells.Item(428,4).Value = 44
$ws.Cells.Item(429,1).Value = '17:14:54'
$ws.Cells.Item(429,2).Value = '18:58'
$ws.Cells.Item(429,3).Value = '215A_EL PATO'
$ws.Cells.Item(429,4).Value = 104
$ws.Cells.Item(430,1).Value = '17:14:54'
$ws.Cells.Item(430,2).Value = '19:04'
$ws.Cells.Item(430,3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(430,4).Value = 110
$ws.Cells.Item(431,1).Value = '18:12:36'
$ws.Cells.Item(431,2).Value = '19:04'
$ws.Cells.Item(431,3).Value = '23_HERNANDEZ'
$ws.Cells.Item(431,4).Value = 52
$ws.Cells.Item(432,1).Value = '17:14:54'
$ws.Cells.Item(432,2).Value = '19:10'
$ws.Cells.Item(432,4).Value = 116
$ws.Cells.Item(433,1).Value = '17:39:57'
$ws.Cells.Item(433,2).Value = '19:16'
$ws.Cells.Item(433,3).Value = '27_EL RETIRO'
$ws.Cells.Item(433,4).Value = 97
$ws.Cells.Item(434,1).Value = '17:39:57'
$ws.Cells.Item(434,2).Value = '19:20'
$ws.Cells.Item(434,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(434,4).Value = 101
$ws.Cells.Item(435,1).Value = '17:53:46'
$ws.Cells.Item(435,2).Value = '19:21'
$ws.Cells.Item(435,3).Value = '26_HERNANDEZ'
$ws.Cells.Item(435,4).Value = 88
$ws.Cells.Item(435,5).Value = 'LP1912'
$ws.Cells.Item(436,1).Value = '17:39:57'
$ws.Cells.Item(436,2).Value = '19:29'
$ws.Cells.Item(436,3).Value = '225_GOMEZ'
$ws.Cells.Item(436,4).Value = 110
$ws.Cells.Item(436,5).Value = 'LP1912'
$ws.Cells.Item(437,1).Value = '17:53:46'
$ws.Cells.Item(437,2).Value = '19:30'
$ws.Cells.Item(437,3).Value = '225_GOMEZ'
$ws.Cells.Item(437,4).Value = 97
$ws.Cells.Item(437,5).Value = 'LP1912'
$ws.Cells.Item(438,1).Value = '17:53:46'
$ws.Cells.Item(438,2).Value = '19:39'
$ws.Cells.Item(438,3).Value = '215C_EL PATO'
$ws.Cells.Item(438,4).Value = 106
$ws.Cells.Item(438,5).Value = 'LP1912'
$ws.Cells.Item(439,1).Value = '17:53:46'
$ws.Cells.Item(439,2).Value = '19:50'
$ws.Cells.Item(439,3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(439,4).Value = 117
$ws.Cells.Item(439,5).Value = 'LP1912'
$ws.Cells.Item(440,1).Value = '17:53:46'
$ws.Cells.Item(440,2).Value = '19:50'
$ws.Cells.Item(440,3).Value = '11X44_ETCHEVERRY'
$ws.Cells.Item(440,4).Value = 117
$ws.Cells.Item(440,5).Value = 'LP1912'
$ws.Cells.Item(441,1).Value = '17:53:46'
$ws.Cells.Item(441,2).Value = '19:51'
$ws.Cells.Item(441,3).Value = '81_EL PELIGRO'
$ws.Cells.Item(441,4).Value = 118
$ws.Cells.Item(441,5).Value = 'LP1912'
$ws.Cells.Item(442,1).Value = '18:12:36'
$ws.Cells.Item(442,2).Value = '19:59'
$ws.Cells.Item(442,3).Value = '17_ROMERO'
$ws.Cells.Item(442,4).Value = 107
$ws.Cells.Item(442,5).Value = 'LP1912'
$ws.Cells.Item(443,1).Value = '18:12:36'
$ws.Cells.Item(443,2).Value = '20:10'
$ws.Cells.Item(443,3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(443,4).Value = 118
$ws.Cells.Item(443,5).Value = 'LP1912'

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2,1).Value = 'Última actualización: 18:12:36'

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2,1).Value = 'Última actualización: 18:12:36'
$ws.Cells.Item(3,1).Value = 'Total filas: 57'
$ws.Cells.Item(59,1).Value = '18:12:36'
$ws.Cells.Item(59,2).Value = '18:52'
$ws.Cells.Item(59,3).Value = '215A_LA PLATA'
$ws.Cells.Item(59,4).Value = 40
$ws.Cells.Item(60,1).Value = '17:14:54'
$ws.Cells.Item(60,2).Value = '19:03'
$ws.Cells.Item(60,3).Value = '215B_LP-P MOR-1 Y 57'
$ws.Cells.Item(60,4).Value = 109
$ws.Cells.Item(60,5).Value = 'L6173'
$ws.Cells.Item(61,1).Value = '18:12:36'
$ws.Cells.Item(61,2).Value = '19:04'
$ws.Cells.Item(61,3).Value = '215B_LP-P MOR-1 Y 57'
$ws.Cells.Item(61,4).Value = 52
$ws.Cells.Item(61,5).Value = 'L6173'
$ws.Cells.Item(62,1).Value = '18:12:36'
$ws.Cells.Item(62,2).Value = '19:53'
$ws.Cells.Item(62,3).Value = '215C_LA PLATA'
$ws.Cells.Item(62,4).Value = 101
$ws.Cells.Item(62,5).Value = 'L6203'
